# Auto-generated edit script applying diff changes to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.635.32'
$ws.Range('E2').Value = '  -0.39%  '

$ws.Range('D3').Value = '2.547.82'
$ws.Range('E3').Value = '  +0.10%  '

$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '317.80'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +4.54%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '95.38'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.63%  '

$ws.Range('E7').Value = '  +0.37%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('E9').Value = '  -1.66%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '36.42'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -0.05%  '

$ws.Range('E11').Value = '  -1.43%  '

$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '7.69'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +1.51%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.114'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.35%  '

$ws.Range('D14').Value = '2.936.09'
$ws.Range('E14').Value = '  -0.08%  '

$ws.Range('E15').Value = '  +4.03%  '

$ws.Range('D16').Value = '2.543.71'
$ws.Range('E16').Value = '  -1.21%  '

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.871'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.31%  '

$ws.Range('D18').Value = '42.702.44'
$ws.Range('E18').Value = '  -0.28%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '13.09'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -1.57%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.67'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +1.24%  '

$ws.Range('D21').Value = '0.0₃0971'
$ws.Range('E21').Value = '  -1.65%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '71.10'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.02%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '256.53'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.44%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.99'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +1.44%  '

$ws.Range('E25').Value = '  -1.42%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '27.57'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -1.71%  '

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.12%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.35'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +3.33%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '39.48'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +3.77%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '10.26'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.66%  '

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '6.00'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -1.88%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '155.97'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.92%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '3.42'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +3.18%  '

$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '2.17'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +1.10%  '

$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '19.39'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -1.25%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.0793'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -0.80%  '

$ws.Range('E37').Value = '  -0.42%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.111'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -4.11%  '

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '24.15'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -5.02%  '

$ws.Range('E40').Value = '  -0.33%  '

$ws.Range('E41').Value = '  +8.15%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '3.85'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.66%  '

$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '3.37'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -1.67%  '

$ws.Range('E44').Value = '  -0.56%  '

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.33%  '

$ws.Range('D46').Value = '2.044.86'
$ws.Range('E46').Value = '  -2.28%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '85.07'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -2.44%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '8.94'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.34%  '

$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.791.04'
$ws.Range('E49').Value = '  -0.13%  '

$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '74.68'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +0.15%  '

$ws.Range('E51').Value = '  -0.41%  '
